# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet (3rd sheet, after Input/Summary) gets a new,
# blank column inserted before what used to be column N ("Late"), pushing the
# old N/O/P ("Late"/"heading"/"Outstanding") columns one to the right
# (-> O/P/Q). The new column takes on the column width that column M
# ("In Advance") already had. The active sheet/selection also moves from the
# "Transactions" sheet to the "Repayment schedule" sheet, with the new
# selection sitting at I16.

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column at column N - this shifts the existing N, O, P
# columns (and all their cell data/styles) one column to the right, becoming
# O, P, Q respectively.
$wsRepayment.Columns.Item(14).Insert()

# The newly inserted column inherits the width of the column immediately to
# its left (column M / "In Advance").
$wsRepayment.Columns.Item(14).ColumnWidth = $wsRepayment.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet (this also clears the
# "tabSelected" flag that used to sit on the "Transactions" sheet, and moves
# the workbook's activeTab accordingly).
$wsRepayment.Activate()

# Restore/update the selected cell on the "Repayment schedule" sheet.
$null = $wsRepayment.Range("I16").Select()
